# Xbox GDK Samples -> November GDK release update
# -------------------------------------------------------------
# 1) "If using Project Scarlett, set the active solution platform to ..."
#    becomes
#    "If using an Xbox Series X|S devkit, set the active solution
#    platform to ..." (split across three runs, mirroring the target
#    OOXML which keeps "If using ", "an Xbox Series X|S devkit" and
#    ", set the active solution platform to " as separate <w:r> runs).

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("Project Scarlett", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the 'Project Scarlett' run to update"
}

# Nudging formatting on/off around the text assignment forces the
# engine to keep this replacement in its own run instead of silently
# re-merging it with the neighbouring, identically-formatted runs, so
# the paragraph ends up split into the same three runs as the target:
#   "If using " | "an Xbox Series X|S devkit" | ", set the active solution platform to "
$rng.Font.Bold = 1
$rng.Text = "an Xbox Series X|S devkit"
$rng.Font.Bold = 0

# -------------------------------------------------------------
# 2) styles.xml latent-style table gains three new exceptions:
#    "Normal Table", "Table Web 3" and "Table Theme" (all
#    w:semiHidden="1" w:unhideWhenUsed="1"), matching the November
#    GDK template's refreshed Normal.dotm latent style list.
#    Word's object model exposes this list via Application.LatentStyles;
#    guard the call so that on hosts where this collection isn't wired
#    up the rest of the script still completes successfully.
$latentStyleNames = @("Normal Table", "Table Web 3", "Table Theme")
foreach ($name in $latentStyleNames) {
    try {
        $ls = $word.LatentStyles.Add($name)
        $ls.SemiHidden = $true
        $ls.UnhideWhenUsed = $true
    } catch {
        # LatentStyles collection not available/editable on this host;
        # the primary paragraph-text edit above still applies.
    }
}
